$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows for QS1 and QS2 (Qiagen 24plex) at the bottom of the table
$ws.Range("A55").Value = "QS1"
$ws.Range("B55").Value = 1
$ws.Range("C55").Value = 1
$ws.Range("D55").Value = 1
$ws.Range("E55").Value = 1
$ws.Range("F55").Value = 1

$ws.Range("A56").Value = "QS2"
$ws.Range("B56").Value = 2
$ws.Range("C56").Value = 2
$ws.Range("D56").Value = 2
$ws.Range("E56").Value = 2
$ws.Range("F56").Value = 2

# Update the view to match the scrolled/selected state in the diff
$excel.ActiveWindow.ScrollRow = 35
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F60").Select()
